{"js": "// Auto-generated body content definition (kept in sync with edit.ps1 via generate_scripts.py)\nconst PARAGRAPHS = [\n  { style: \"Heading1\", runs: [\n      { text: \"Kn\u00e4rot \u2013 ekologi samt krav p\u00e5 livsmilj\u00f6n\", italic: false }\n  ] },\n  { style: null, runs: [\n      { text: \"Kn\u00e4rot \u00e4r fridlyst enligt 8 och 15 \u00a7\u00a7 artskyddsf\u00f6rordningen och klassad som s\u00e5rbar (VU) enligt r\u00f6dlistan 2020. Kn\u00e4rot \u00e4r beroende av h\u00f6g och j\u00e4mn luftfuktighet i gamla, ost\u00f6rda skogsmilj\u00f6er och \u00e4r k\u00e4nslig f\u00f6r snabba f\u00f6r\u00e4ndringar av ljus-/vindf\u00f6rh\u00e5llanden eller uttorkning. P\u00e5 grund av ett alltf\u00f6r intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 \u00e5ren och i framtiden bed\u00f6ms minskningstakten uppg\u00e5 till 30 (20-40) %. Till f\u00f6ljd av att arten har en dokumenterat h\u00f6gre minskningstakt if\u00f6rh\u00e5llande till sin generationstid \u00e4n vad som tidigare varit k\u00e4nt (data fr\u00e5n Riksskogstaxeringen) h\u00f6jdes den till hotkategori s\u00e5rbar (VU) i r\u00f6dlistan 2020 (Artdatabanken, 2021).\", italic: false }\n  ] },\n  { style: null, runs: [\n      { text: \"Samuel Johnsons doktorsavhandling \", italic: false },\n      { text: \"\u201cRetention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation\u201c\", italic: true },\n      { text: \" (SLU, Uppsala 2014) visar att det kr\u00e4vs v\u00e4l tilltagna skyddszoner f\u00f6r att kn\u00e4rotens v\u00e4xtplatser inte ska ta skada av skogsbruks\u00e5tg\u00e4rder i intilliggande omr\u00e5den: \", italic: false },\n      { text: \"\u201cStudy III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.\u201d \", italic: true },\n      { text: \"Vidare \", italic: false },\n      { text: \"\u201cMore sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).\u201d\", italic: true }\n  ] },\n  { style: null, runs: [\n      { text: \"Johnsons (2014) rekommendation p\u00e5 minst 50 meters breda skyddszoner runt kn\u00e4rotens v\u00e4xtplatser motsvarar en areal p\u00e5 0,78 hektar, vilket ligger i linje med andra studier som gjorts p\u00e5 k\u00e4nsliga skogsarter: \", italic: false },\n      { text: \"\u201cIn study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).\u201d\", italic: true }\n  ] },\n  { style: null, runs: [\n      { text: \"En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkid\u00e9n kn\u00e4rots skyddsbehov. I uppsatsen ber\u00f6rs problemet med uttorkning f\u00f6r v\u00e4xter, bl.a. f\u00f6r kn\u00e4rot, ett problem som blivit accentuerat p\u00e5 grund av den p\u00e5g\u00e5ende klimatf\u00f6r\u00e4ndringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen unders\u00f6ks omr\u00e5den med tre olika avst\u00e5nd fr\u00e5n kalhyggeskant med avseende p\u00e5 skydd bl.a. f\u00f6r kn\u00e4rot. Det f\u00f6rsta omr\u00e5det har avst\u00e5nd upp till 20 m fr\u00e5n hyggeskant (Strong edge effect), det andra 20 \u2013 40 m fr\u00e5n hyggeskant (Weak edge effect) och det tredje avser st\u00f6rre avst\u00e5nd fr\u00e5n hyggeskant, d\u00e4r kanteffekten anses vara f\u00f6rsumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt p\u00e5 k\u00e4nsliga och r\u00f6dlistade skogsarter vid de kortare avst\u00e5nden till hyggeskant, medan effekt av uttorkning inte konstaterades p\u00e5 st\u00f6rre avst\u00e5nd (Interior). F\u00f6r orkid\u00e9n kn\u00e4rot fann man en rik f\u00f6rekomst (upp till 0,06 dm2/m2) p\u00e5 stort avst\u00e5nd fr\u00e5n hyggeskant (Interior), medan f\u00f6rekomsten var liten eller n\u00e4rmast f\u00f6rsumbar i de omr\u00e5den som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet p\u00e5pekar att de allt oftare f\u00f6rekommande torra somrarna ger ytterligare sk\u00e4l att ut\u00f6ka skyddsavst\u00e5ndet fr\u00e5n hyggen till den fuktkr\u00e4vande arten kn\u00e4rot (Koelmeijer m.fl., 2022).\", italic: false }\n  ] },\n  { style: null, runs: [\n      { text: \"\u00c4ven Skogsstyrelsens egen v\u00e4gledning f\u00f6r h\u00e4nsyn till kn\u00e4rot ligger i linje med ovanst\u00e5ende forskningsstudier. Av v\u00e4gledningen framg\u00e5r det att f\u00f6r med h\u00f6g sannolikhet kunna bevara befintliga f\u00f6rekomster kr\u00e4vs relativt stora avs\u00e4ttningar av uppvuxen skog med slutet och relativt t\u00e4tt kronskikt. Som riktlinje kan kr\u00e4vas ett avst\u00e5nd p\u00e5 50 meter in fr\u00e5n brynet f\u00f6r att vidmakth\u00e5lla ett fungerande mikroklimat. Detta inneb\u00e4r att frist\u00e5ende h\u00e4nsynsytor f\u00f6r m\u00e5nga arter (k\u00e4rlv\u00e4xter, lavar och mossor) kan beh\u00f6va ha en area \u00f6verstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) f\u00f6r att bibeh\u00e5lla lokalklimatet. \u00c4ven ganska sm\u00e5 f\u00f6r\u00e4ndringar i form av f\u00f6r\u00e4ndrade ljus- och fuktighetsf\u00f6rh\u00e5llanden, till exempel till f\u00f6ljd av gallring, kan leda till att arten f\u00f6rsvinner till f\u00f6ljd av konkurrens med mera ljuskr\u00e4vande och snabbv\u00e4xande arter (Skogsstyrelsen, 2022).\", italic: false }\n  ] },\n  { style: \"Heading2\", runs: [\n      { text: \"Referenser - kn\u00e4rot\", italic: false }\n  ] },\n  { style: null, runs: [\n      { text: \"de Graaf M & Roberts M.R., 2009. \", italic: false },\n      { text: \"Short-term response of the herbaceous layer within leave patches after harvest. \", italic: true },\n      { text: \"Forest Ecology and Management 257, 1014-1025\", italic: false }\n  ] },\n  { style: null, runs: [\n      { text: \"Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. \", italic: false },\n      { text: \"Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. \", italic: true },\n      { text: \"Ecological Applications, 22, 2049-2064 \", italic: false }\n  ] },\n  { style: null, runs: [\n      { text: \"Koelmeijer, I. A., Ehrl\u00e9n, J., J\u00f6nsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. \", italic: false },\n      { text: \"Interactive effects of drought and edge exposure on old-growth forest understory species. \", italic: true },\n      { text: \"Landscape Ecology, 37, sid 1839-1853\", italic: false }\n  ] },\n  { style: null, runs: [\n      { text: \"Rudolphi, J., J\u00f6nsson, M. T., & Gustafsson, L., 2014. \", italic: false },\n      { text: \"Biological legacies buffer local species extinction after logging. \", italic: true },\n      { text: \"Journal of Applied Ecology. 51, 53-62.\", italic: false }\n  ] },\n  { style: null, runs: [\n      { text: \"Skogsstyrelsen, 2022. \", italic: false },\n      { text: \"V\u00e4gledning f\u00f6r h\u00e4nsyn till kn\u00e4rot. \", italic: true },\n      { text: \"https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/\", italic: false }\n  ] },\n  { style: null, runs: [\n      { text: \"SLU Artdatabanken, 2021. \", italic: false },\n      { text: \"Artfaktablad. Naturv\u00e5rd \u2013 artfakta. \", italic: true },\n      { text: \"SLU Artdatabanken, Uppsala \", italic: false }\n  ] },\n];\n\n// 1. Insert the new \"kn\u00e4rot\" section right after the \"BILAGA 1 - Fridlysta arter\"\n//    paragraph, and before the final section break.\nconst bodyParagraphs = context.document.body.paragraphs;\nbodyParagraphs.load(\"items/text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = 0; i < bodyParagraphs.items.length; i++) {\n  if (bodyParagraphs.items[i].text === \"BILAGA 1 - Fridlysta arter\") {\n    anchor = bodyParagraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not find anchor paragraph \"BILAGA 1 - Fridlysta arter\"');\n}\n\nlet insertAfter = anchor;\nfor (const para of PARAGRAPHS) {\n  // Create an empty paragraph right after the current insertion point.\n  const newPara = insertAfter.insertParagraph(\"\", Word.InsertLocation.after);\n  // Explicitly (re)set the style: either a heading style, or \"Normal\" so that\n  // no pPr/pStyle is inherited from the previous paragraph (e.g. \"Title\").\n  newPara.style = para.style ? para.style : \"Normal\";\n\n  // Add each run of text, toggling italics as required. We insert each run\n  // at the end of the (so far empty) paragraph, in order.\n  for (const run of para.runs) {\n    const insertedRange = newPara.insertText(run.text, Word.InsertLocation.end);\n    if (run.italic) {\n      insertedRange.font.italic = true;\n    }\n  }\n\n  insertAfter = newPara;\n}\n\nawait context.sync();\n\n// 2. Update the date in the first-page header from 2023-09-13 to 2023-09-15.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let s = 0; s < sections.items.length; s++) {\n  const header = sections.items[s].getHeader(Word.HeaderFooterType.firstPage);\n  const results = header.search(\"2023-09-13\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let r = 0; r < results.items.length; r++) {\n    results.items[r].insertText(\"2023-09-15\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Auto-generated body content definition (kept in sync with edit.js via generate_scripts.py)\n$paragraphs = @(\n    @{\n        Style = \"Heading1\"\n        Runs = @(\n        @{ Text = \"Kn\u00e4rot \u2013 ekologi samt krav p\u00e5 livsmilj\u00f6n\"; Italic = $false }\n        )\n    },\n    @{\n        Style = $null\n        Runs = @(\n        @{ Text = \"Kn\u00e4rot \u00e4r fridlyst enligt 8 och 15 \u00a7\u00a7 artskyddsf\u00f6rordningen och klassad som s\u00e5rbar (VU) enligt r\u00f6dlistan 2020. Kn\u00e4rot \u00e4r beroende av h\u00f6g och j\u00e4mn luftfuktighet i gamla, ost\u00f6rda skogsmilj\u00f6er och \u00e4r k\u00e4nslig f\u00f6r snabba f\u00f6r\u00e4ndringar av ljus-/vindf\u00f6rh\u00e5llanden eller uttorkning. P\u00e5 grund av ett alltf\u00f6r intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 \u00e5ren och i framtiden bed\u00f6ms minskningstakten uppg\u00e5 till 30 (20-40) %. Till f\u00f6ljd av att arten har en dokumenterat h\u00f6gre minskningstakt if\u00f6rh\u00e5llande till sin generationstid \u00e4n vad som tidigare varit k\u00e4nt (data fr\u00e5n Riksskogstaxeringen) h\u00f6jdes den till hotkategori s\u00e5rbar (VU) i r\u00f6dlistan 2020 (Artdatabanken, 2021).\"; Italic = $false }\n        )\n    },\n    @{\n        Style = $null\n        Runs = @(\n        @{ Text = \"Samuel Johnsons doktorsavhandling \"; Italic = $false },\n        @{ Text = \"\u201cRetention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation\u201c\"; Italic = $true },\n        @{ Text = \" (SLU, Uppsala 2014) visar att det kr\u00e4vs v\u00e4l tilltagna skyddszoner f\u00f6r att kn\u00e4rotens v\u00e4xtplatser inte ska ta skada av skogsbruks\u00e5tg\u00e4rder i intilliggande omr\u00e5den: \"; Italic = $false },\n        @{ Text = \"\u201cStudy III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.\u201d \"; Italic = $true },\n        @{ Text = \"Vidare \"; Italic = $false },\n        @{ Text = \"\u201cMore sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).\u201d\"; Italic = $true }\n        )\n    },\n    @{\n        Style = $null\n        Runs = @(\n        @{ Text = \"Johnsons (2014) rekommendation p\u00e5 minst 50 meters breda skyddszoner runt kn\u00e4rotens v\u00e4xtplatser motsvarar en areal p\u00e5 0,78 hektar, vilket ligger i linje med andra studier som gjorts p\u00e5 k\u00e4nsliga skogsarter: \"; Italic = $false },\n        @{ Text = \"\u201cIn study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).\u201d\"; Italic = $true }\n        )\n    },\n    @{\n        Style = $null\n        Runs = @(\n        @{ Text = \"En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkid\u00e9n kn\u00e4rots skyddsbehov. I uppsatsen ber\u00f6rs problemet med uttorkning f\u00f6r v\u00e4xter, bl.a. f\u00f6r kn\u00e4rot, ett problem som blivit accentuerat p\u00e5 grund av den p\u00e5g\u00e5ende klimatf\u00f6r\u00e4ndringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen unders\u00f6ks omr\u00e5den med tre olika avst\u00e5nd fr\u00e5n kalhyggeskant med avseende p\u00e5 skydd bl.a. f\u00f6r kn\u00e4rot. Det f\u00f6rsta omr\u00e5det har avst\u00e5nd upp till 20 m fr\u00e5n hyggeskant (Strong edge effect), det andra 20 \u2013 40 m fr\u00e5n hyggeskant (Weak edge effect) och det tredje avser st\u00f6rre avst\u00e5nd fr\u00e5n hyggeskant, d\u00e4r kanteffekten anses vara f\u00f6rsumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt p\u00e5 k\u00e4nsliga och r\u00f6dlistade skogsarter vid de kortare avst\u00e5nden till hyggeskant, medan effekt av uttorkning inte konstaterades p\u00e5 st\u00f6rre avst\u00e5nd (Interior). F\u00f6r orkid\u00e9n kn\u00e4rot fann man en rik f\u00f6rekomst (upp till 0,06 dm2/m2) p\u00e5 stort avst\u00e5nd fr\u00e5n hyggeskant (Interior), medan f\u00f6rekomsten var liten eller n\u00e4rmast f\u00f6rsumbar i de omr\u00e5den som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet p\u00e5pekar att de allt oftare f\u00f6rekommande torra somrarna ger ytterligare sk\u00e4l att ut\u00f6ka skyddsavst\u00e5ndet fr\u00e5n hyggen till den fuktkr\u00e4vande arten kn\u00e4rot (Koelmeijer m.fl., 2022).\"; Italic = $false }\n        )\n    },\n    @{\n        Style = $null\n        Runs = @(\n        @{ Text = \"\u00c4ven Skogsstyrelsens egen v\u00e4gledning f\u00f6r h\u00e4nsyn till kn\u00e4rot ligger i linje med ovanst\u00e5ende forskningsstudier. Av v\u00e4gledningen framg\u00e5r det att f\u00f6r med h\u00f6g sannolikhet kunna bevara befintliga f\u00f6rekomster kr\u00e4vs relativt stora avs\u00e4ttningar av uppvuxen skog med slutet och relativt t\u00e4tt kronskikt. Som riktlinje kan kr\u00e4vas ett avst\u00e5nd p\u00e5 50 meter in fr\u00e5n brynet f\u00f6r att vidmakth\u00e5lla ett fungerande mikroklimat. Detta inneb\u00e4r att frist\u00e5ende h\u00e4nsynsytor f\u00f6r m\u00e5nga arter (k\u00e4rlv\u00e4xter, lavar och mossor) kan beh\u00f6va ha en area \u00f6verstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) f\u00f6r att bibeh\u00e5lla lokalklimatet. \u00c4ven ganska sm\u00e5 f\u00f6r\u00e4ndringar i form av f\u00f6r\u00e4ndrade ljus- och fuktighetsf\u00f6rh\u00e5llanden, till exempel till f\u00f6ljd av gallring, kan leda till att arten f\u00f6rsvinner till f\u00f6ljd av konkurrens med mera ljuskr\u00e4vande och snabbv\u00e4xande arter (Skogsstyrelsen, 2022).\"; Italic = $false }\n        )\n    },\n    @{\n        Style = \"Heading2\"\n        Runs = @(\n        @{ Text = \"Referenser - kn\u00e4rot\"; Italic = $false }\n        )\n    },\n    @{\n        Style = $null\n        Runs = @(\n        @{ Text = \"de Graaf M & Roberts M.R., 2009. \"; Italic = $false },\n        @{ Text = \"Short-term response of the herbaceous layer within leave patches after harvest. \"; Italic = $true },\n        @{ Text = \"Forest Ecology and Management 257, 1014-1025\"; Italic = $false }\n        )\n    },\n    @{\n        Style = $null\n        Runs = @(\n        @{ Text = \"Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. \"; Italic = $false },\n        @{ Text = \"Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. \"; Italic = $true },\n        @{ Text = \"Ecological Applications, 22, 2049-2064 \"; Italic = $false }\n        )\n    },\n    @{\n        Style = $null\n        Runs = @(\n        @{ Text = \"Koelmeijer, I. A., Ehrl\u00e9n, J., J\u00f6nsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. \"; Italic = $false },\n        @{ Text = \"Interactive effects of drought and edge exposure on old-growth forest understory species. \"; Italic = $true },\n        @{ Text = \"Landscape Ecology, 37, sid 1839-1853\"; Italic = $false }\n        )\n    },\n    @{\n        Style = $null\n        Runs = @(\n        @{ Text = \"Rudolphi, J., J\u00f6nsson, M. T., & Gustafsson, L., 2014. \"; Italic = $false },\n        @{ Text = \"Biological legacies buffer local species extinction after logging. \"; Italic = $true },\n        @{ Text = \"Journal of Applied Ecology. 51, 53-62.\"; Italic = $false }\n        )\n    },\n    @{\n        Style = $null\n        Runs = @(\n        @{ Text = \"Skogsstyrelsen, 2022. \"; Italic = $false },\n        @{ Text = \"V\u00e4gledning f\u00f6r h\u00e4nsyn till kn\u00e4rot. \"; Italic = $true },\n        @{ Text = \"https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/\"; Italic = $false }\n        )\n    },\n    @{\n        Style = $null\n        Runs = @(\n        @{ Text = \"SLU Artdatabanken, 2021. \"; Italic = $false },\n        @{ Text = \"Artfaktablad. Naturv\u00e5rd \u2013 artfakta. \"; Italic = $true },\n        @{ Text = \"SLU Artdatabanken, Uppsala \"; Italic = $false }\n        )\n    },\n)\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# 1. Insert the new \"kn\u00e4rot\" section right after the \"BILAGA 1 - Fridlysta\n#    arter\" paragraph (the last paragraph of the body), and before the final\n#    section break.\n# ---------------------------------------------------------------------------\n$anchor = $null\nforeach ($para in $d.Paragraphs) {\n    $t = $para.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"BILAGA 1 - Fridlysta arter\") {\n        $anchor = $para\n        break\n    }\n}\n\nif ($null -eq $anchor) {\n    throw 'Could not find anchor paragraph \"BILAGA 1 - Fridlysta arter\"'\n}\n\n# $r is a \"live\" insertion-point range that we keep collapsed at the end of\n# the most-recently-inserted content; every new paragraph is appended there.\n$r = $anchor.Range.Duplicate\n$r.Collapse(0)  # wdCollapseEnd\n\nforeach ($paraDef in $paragraphs) {\n    # Create a new (empty) paragraph right after the current insertion point.\n    $r.InsertParagraphAfter()\n    $r.Collapse(0)\n    $r.Move(4, 1) | Out-Null  # wdParagraph: step into the paragraph just created\n\n    $newPara = $d.Paragraphs.Last\n\n    # Explicitly (re)apply the paragraph style via ParagraphFormat so that a\n    # \"Normal\" paragraph has no leftover pPr/pStyle (matching a plain <w:p>),\n    # instead of inheriting the previous (e.g. Title/Heading) style.\n    if ($paraDef.Style) {\n        $newPara.Format.Style = $paraDef.Style\n    } else {\n        $newPara.Format.Style = \"Normal\"\n    }\n\n    # Insert the whole paragraph's text as one chunk (concatenation of all\n    # runs), then go back and italicize the relevant substrings using Find\n    # scoped to this paragraph's range. Doing it this way (rather than\n    # toggling Font.Italic on a live collapsed cursor) avoids Word's\n    # \"formatting bleeds into the next paragraph/run\" quirk.\n    $fullText = ($paraDef.Runs | ForEach-Object { $_.Text }) -join \"\"\n    $r.InsertAfter($fullText)\n    $r.Collapse(0)\n\n    $paraRange = $newPara.Range\n    $cursor = $paraRange.Duplicate\n    $cursor.Collapse(1)  # wdCollapseStart\n\n    foreach ($run in $paraDef.Runs) {\n        $searchRange = $d.Range($cursor.Start, $paraRange.End)\n        $searchRange.Find.ClearFormatting()\n        $searchRange.Find.MatchCase = $true\n        $searchRange.Find.MatchWildcards = $false\n        $searchRange.Find.Text = $run.Text\n        $found = $searchRange.Find.Execute()\n        if (-not $found) {\n            throw (\"Run text not found while formatting: \" + $run.Text)\n        }\n        if ($run.Italic) {\n            $searchRange.Font.Italic = 1\n        }\n        $cursor = $searchRange.Duplicate\n        $cursor.Collapse(0)  # wdCollapseEnd: continue searching after this run\n    }\n}\n\n# ---------------------------------------------------------------------------\n# 2. Update the date in the first-page header from 2023-09-13 to 2023-09-15.\n# ---------------------------------------------------------------------------\nforeach ($sec in $d.Sections) {\n    $hdr = $sec.Headers.Item(2)  # wdHeaderFooterFirstPage\n    $find = $hdr.Range.Find\n    $find.ClearFormatting()\n    $find.Text = \"2023-09-13\"\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = \"2023-09-15\"\n    $find.Execute($null, $true, $true, $false, $null, $null, $true, $null, $null, $null, 2) | Out-Null  # wdReplaceAll\n}\n"}
